$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before column B - this shifts the existing
#    B:E data (am/working/with/openpyxl) right into C:F.
$ws.Columns("B:B").Insert()

# 2) Clear out rows 9:13 (their content, leaving the empty rows in place).
$ws.Range("A9:F13").ClearContents()
$ws.Rows("9:13").OutlineLevel = 0

# 3) Append new rows 21:24 with the same "I am working with openpyxl"
#    pattern (column B intentionally left blank, matching the rest of
#    the sheet after the column insert).
for ($r = 21; $r -le 24; $r++) {
    $ws.Cells.Item($r, 1).Value = "I"
    $ws.Cells.Item($r, 3).Value = "am"
    $ws.Cells.Item($r, 4).Value = "working"
    $ws.Cells.Item($r, 5).Value = "with"
    $ws.Cells.Item($r, 6).Value = "openpyxl"
}
